$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (group by date processing script changes)
$ws.Range("F2").Value = 0
$ws.Range("J2").Value = 1

# Update the view: scroll so column D is the top-left visible column,
# and set the active selection to J3
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("J3").Select()
